$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.005.22"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -2.35%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.102.00"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -1.10%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -1.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "346.84"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +2.45%  "
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -1.00%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5156"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -2.40%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4430"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -3.33%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.09345"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +2.39%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "52.57"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -3.75%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.174"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -0.38%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "25.30"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +2.92%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.097.51"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -1.66%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.750"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -1.86%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "8.184"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +0.47%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "99.53"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +2.17%  "
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -1.63%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.004"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -1.08%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "20.63"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +5.40%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.06680"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -0.47%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.003"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -1.05%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.229"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -3.23%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "30.079.71"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -2.38%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.66"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -2.77%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.333"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -1.23%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.345.06"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -1.47%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "22.06"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -2.38%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.555"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -0.34%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "162.84"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -1.65%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "133.56"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -1.08%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.169"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -3.54%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.1059"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -1.97%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.643"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -0.92%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.242"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -2.35%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.943"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.192"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +4.56%  "
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -4.43%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02568"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -4.51%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06782"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -1.57%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.2283"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -2.28%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "12.59"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -0.77%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.6926"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -0.20%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.310"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +3.56%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.6653"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +2.40%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "14.15"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -6.42%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.292"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -1.15%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.633"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -2.03%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.00000000353"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -6.17%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.222"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -3.31%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "82.13"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -2.17%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.07208"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -1.45%  "
